$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (C) and montant_total (D) columns for rows with new/updated 2020-08-27 data
# Values are kept as text to match the original inlineStr string cell formatting, then
# formats are cleared so no extra style index is left on the cell.
$updates = @(
    @{ Row = 3; C = "1021"; D = "3260894.33" }
    @{ Row = 4; C = "422"; D = "1744198.25" }
    @{ Row = 5; C = "119"; D = "574128.09" }
    @{ Row = 34; C = "575"; D = "1894813.66" }
    @{ Row = 35; C = "230"; D = "1159788.11" }
    @{ Row = 36; C = "74"; D = "407894.00" }
    @{ Row = 52; C = "598"; D = "2117095.21" }
    @{ Row = 54; C = "89"; D = "520878.23" }
    @{ Row = 57; C = "715"; D = "1837318.62" }
    @{ Row = 61; C = "129"; D = "888623.00" }
    @{ Row = 64; C = "16"; D = "41261.00" }
    @{ Row = 82; C = "230"; D = "593326.09" }
    @{ Row = 83; C = "894"; D = "2868012.26" }
    @{ Row = 85; C = "118"; D = "581984.52" }
    @{ Row = 94; C = "97"; D = "254578.00" }
    @{ Row = 95; C = "412"; D = "1250108.52" }
    @{ Row = 96; C = "176"; D = "707710.27" }
    @{ Row = 97; C = "58"; D = "296911.73" }
    @{ Row = 99; C = "14"; D = "28000.00" }
    @{ Row = 100; C = "304"; D = "793054.43" }
    @{ Row = 101; C = "1232"; D = "3782146.74" }
    @{ Row = 102; C = "458"; D = "1883812.62" }
    @{ Row = 103; C = "124"; D = "593996.00" }
    @{ Row = 104; C = "35"; D = "229157.00" }
    @{ Row = 105; C = "67"; D = "146571.16" }
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.C
    $cCell.ClearFormats()

    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.ClearFormats()
}

Write-Host "Updated $($updates.Count) rows"
